# Styling elements dynamically with ngStyle
#
# Split the paragraph right after "...doesn't add it." (just before the
# _GoBack bookmark) into two paragraphs, then add the new explanatory
# sentence about attribute directives into the new paragraph, keeping the
# _GoBack bookmark in place between "on. They" and " look like normal...".

$d = $word.ActiveDocument

# Locate the end of "...it doesn't add it." (right before the bookmark).
$found = $d.Content
$found.Find.Execute("doesn’t add it.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# Split the paragraph there: insert a paragraph mark right before the bookmark.
$splitPoint = $d.Range($found.End, $found.End)
$splitPoint.InsertBefore("`r")

# The bookmark now starts the new (second) paragraph.
$bm = $d.Bookmarks("_GoBack")
$pos = $bm.Start

# Insert the trailing sentence fragment after the bookmark first (position
# is still valid / unaffected by insertions that happen after it).
$afterBookmark = $d.Range($pos, $pos)
$afterBookmark.InsertAfter(" look like normal HTML attributes without a star basically.")

# Insert the text before the bookmark, as two runs (mirrors original authoring).
$beforeBookmark2 = $d.Range($pos, $pos)
$beforeBookmark2.InsertBefore("on. They")
$beforeBookmark1 = $d.Range($pos, $pos)
$beforeBookmark1.InsertBefore("Unlike structural directives, attribute directives don’t add or remove elements. They only change the elements they are placed ")
